$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 6-9 entirely (shifts nothing below them, just drops them)
$ws.Range("A6:B9").EntireRow.Delete()

# Update the surviving data rows (2-5) with their new values
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 278

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = 200

$ws.Range("A4").Value = 21
$ws.Range("B4").Value = 135

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 17
